# Update cryptos list data cells (Price / Volume(1h), and one coin row
# replacement) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.496.64'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.625.84'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.41%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.98'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.45'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.624.48'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.48'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.06%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.365'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.064.89'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.495.56'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.30'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000141'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.612.58'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.40'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +9.34%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '346.61'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.02'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +8.37%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +14.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.29'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.73'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0787'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.83'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +9.54%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.23%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '161.00'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.79%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.07%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.975'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +10.19%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.14%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +7.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.83'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.83'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.849'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '295.89'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '137.86'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0985'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.78%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0545'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.84'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.35%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.72'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.80%  '
